# Generate Report for Handoff
#
# The source file e2e\03568c1b-f54f-4026-b3be-58b205e1bdfa.md was renamed /
# regenerated as e2e\ca991a30-e8df-4831-81ef-906dc091c117.md. A fresh
# handoff (xliff) pair was produced for it, so the report's "latest
# handoff" columns move forward while the "latest target / handback"
# columns reset because nothing has been handed back yet for the new file.

$wb = $excel.ActiveWorkbook

$oldGuid = "03568c1b-f54f-4026-b3be-58b205e1bdfa"
$newGuid = "ca991a30-e8df-4831-81ef-906dc091c117"
$oldCommit = "4a0b53637e967a114459b66d00793d372bc3d38c"
$newCommit = "5a3e5432f83950b9344bd3684f8af4ceef273807"

function Set-DisplayHyperlink {
    param($ws, $addr, $newDisplay)
    foreach ($hl in $ws.Hyperlinks) {
        $a = $hl.Range.Address()
        if ($a -eq $addr) {
            $hl.TextToDisplay = $newDisplay
        }
    }
}

function Remove-HyperlinkAt {
    param($ws, $addr)
    foreach ($hl in $ws.Hyperlinks) {
        $a = $hl.Range.Address()
        if ($a -eq $addr) {
            $hl.Delete()
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = "$newGuid.md"
$ov.Range("B2").Value = "e2e\$newGuid.md"
Set-DisplayHyperlink $ov '$B$2' "e2e\$newGuid.md"
$ov.Range("G2").Value = "2016-08-31 19:08:46"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Value = "$newGuid.md"
Set-DisplayHyperlink $zh '$A$2' "$newGuid.md"

$zh.Range("G2").Value = "$newGuid.$newCommit.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-31 19:08:42"
Remove-HyperlinkAt $zh '$I$2'
$zh.Range("I2").Value = ""
$zh.Range("I2").Style = "Normal"
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

$zh.Columns.Item(9).ColumnWidth = 18.65
$zh.Columns.Item(10).ColumnWidth = 21.71

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Value = "$newGuid.md"
Set-DisplayHyperlink $de '$A$2' "$newGuid.md"

$de.Range("G2").Value = "$newGuid.$newCommit.de-de.xlf"
$de.Range("H2").Value = "2016-08-31 19:08:46"
Remove-HyperlinkAt $de '$I$2'
$de.Range("I2").Value = ""
$de.Range("I2").Style = "Normal"
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Columns.Item(9).ColumnWidth = 18.65
$de.Columns.Item(10).ColumnWidth = 21.71
